{"js": "// Replace the date line and each two-digit-multiplication answer cell with\n// its new value. Every \"old\" value below is unique within the document, so\n// a plain-text body.search() + insertText(\"Replace\") round-trip is safe for\n// each pair (no risk of touching an already-updated cell, since old/new\n// sets don't overlap in this document).\nconst replacements = [\n  [\"2026-02-10 Tuesday\", \"2026-02-11 Wednesday\"],\n  [\"83\u00d727=2241\", \"29\u00d714=406\"],\n  [\"47\u00d719=893\", \"38\u00d744=1672\"],\n  [\"80\u00d719=1520\", \"21\u00d746=966\"],\n  [\"31\u00d746=1426\", \"96\u00d785=8160\"],\n  [\"21\u00d783=1743\", \"75\u00d765=4875\"],\n  [\"49\u00d755=2695\", \"87\u00d751=4437\"],\n  [\"96\u00d745=4320\", \"68\u00d715=1020\"],\n  [\"27\u00d727=729\", \"38\u00d780=3040\"],\n  [\"55\u00d797=5335\", \"43\u00d763=2709\"],\n  [\"47\u00d788=4136\", \"49\u00d719=931\"],\n  [\"51\u00d713=663\", \"51\u00d771=3621\"],\n  [\"40\u00d777=3080\", \"67\u00d763=4221\"],\n  [\"25\u00d714=350\", \"81\u00d777=6237\"],\n  [\"17\u00d786=1462\", \"20\u00d752=1040\"],\n  [\"54\u00d743=2322\", \"38\u00d713=494\"],\n  [\"77\u00d743=3311\", \"22\u00d789=1958\"],\n  [\"63\u00d750=3150\", \"94\u00d730=2820\"],\n  [\"46\u00d762=2852\", \"74\u00d762=4588\"],\n  [\"15\u00d778=1170\", \"27\u00d754=1458\"],\n  [\"63\u00d772=4536\", \"26\u00d770=1820\"],\n  [\"97\u00d760=5820\", \"36\u00d777=2772\"],\n  [\"79\u00d749=3871\", \"84\u00d732=2688\"],\n  [\"95\u00d714=1330\", \"87\u00d773=6351\"],\n  [\"35\u00d731=1085\", \"61\u00d736=2196\"],\n  [\"49\u00d766=3234\", \"29\u00d770=2030\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each two-digit-multiplication answer cell with\n# its new value. Every \"old\" value is unique within the document, so a\n# plain Find/Replace (ReplaceAll) round-trip is safe for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-10 Tuesday\", \"2026-02-11 Wednesday\"),\n    @(\"83\u00d727=2241\", \"29\u00d714=406\"),\n    @(\"47\u00d719=893\", \"38\u00d744=1672\"),\n    @(\"80\u00d719=1520\", \"21\u00d746=966\"),\n    @(\"31\u00d746=1426\", \"96\u00d785=8160\"),\n    @(\"21\u00d783=1743\", \"75\u00d765=4875\"),\n    @(\"49\u00d755=2695\", \"87\u00d751=4437\"),\n    @(\"96\u00d745=4320\", \"68\u00d715=1020\"),\n    @(\"27\u00d727=729\", \"38\u00d780=3040\"),\n    @(\"55\u00d797=5335\", \"43\u00d763=2709\"),\n    @(\"47\u00d788=4136\", \"49\u00d719=931\"),\n    @(\"51\u00d713=663\", \"51\u00d771=3621\"),\n    @(\"40\u00d777=3080\", \"67\u00d763=4221\"),\n    @(\"25\u00d714=350\", \"81\u00d777=6237\"),\n    @(\"17\u00d786=1462\", \"20\u00d752=1040\"),\n    @(\"54\u00d743=2322\", \"38\u00d713=494\"),\n    @(\"77\u00d743=3311\", \"22\u00d789=1958\"),\n    @(\"63\u00d750=3150\", \"94\u00d730=2820\"),\n    @(\"46\u00d762=2852\", \"74\u00d762=4588\"),\n    @(\"15\u00d778=1170\", \"27\u00d754=1458\"),\n    @(\"63\u00d772=4536\", \"26\u00d770=1820\"),\n    @(\"97\u00d760=5820\", \"36\u00d777=2772\"),\n    @(\"79\u00d749=3871\", \"84\u00d732=2688\"),\n    @(\"95\u00d714=1330\", \"87\u00d773=6351\"),\n    @(\"35\u00d731=1085\", \"61\u00d736=2196\"),\n    @(\"49\u00d766=3234\", \"29\u00d770=2030\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Saved = $false\n"}
